# Simulated Wild Card round and logged it:
# update the Rushing and Receiving stat sheets with the round's numbers.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# M.Ryan (row 2)
$rushing.Range("C2").Value = 4

# M.Davis (row 4)
$rushing.Range("C4").Value = 92
$rushing.Range("D4").Value = 37
$rushing.Range("E4").Value = 9
$rushing.Range("F4").Value = 20

# C.Patterson (row 5)
$rushing.Range("C5").Value = 103
$rushing.Range("D5").Value = 39
$rushing.Range("E5").Value = 12

# Q.Ollison (row 7)
$rushing.Range("D7").Value = 4
$rushing.Range("F7").Value = 3

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# M.Davis (row 2)
$receiving.Range("C2").Value = 56
$receiving.Range("D2").Value = 42

# C.Patterson (row 3)
$receiving.Range("C3").Value = 60
$receiving.Range("D3").Value = 45

# Q.Ollison (row 5)
$receiving.Range("C5").Value = 6
$receiving.Range("G5").Value = 2

# R.Gage (row 6)
$receiving.Range("C6").Value = 83
$receiving.Range("D6").Value = 71
$receiving.Range("E6").Value = 22
$receiving.Range("F6").Value = 12
$receiving.Range("G6").Value = 14
$receiving.Range("H6").Value = 9

# O.Zaccheaus (row 7)
$receiving.Range("C7").Value = 44
$receiving.Range("D7").Value = 26
$receiving.Range("E7").Value = 7
$receiving.Range("F7").Value = 4

# C.Blake (row 8)
$receiving.Range("C8").Value = 7

# T.Sharpe (row 9)
$receiving.Range("C9").Value = 30
$receiving.Range("E9").Value = 6
$receiving.Range("G9").Value = 5

# K.Pitts (row 12)
$receiving.Range("C12").Value = 80
$receiving.Range("D12").Value = 52
$receiving.Range("E12").Value = 29
$receiving.Range("G12").Value = 15

# H.Hurst (row 13)
$receiving.Range("E13").Value = 2
$receiving.Range("F13").Value = 1

# K.Smith (row 16)
$receiving.Range("C16").Value = 6
$receiving.Range("D16").Value = 6
